$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Note: literal text values that look like dates ("DD.MM.YYYY") get silently
# auto-converted into date serial numbers if assigned directly via .Value.
# To avoid that, we put the literal text behind a formula (so it is typed as
# text from the start), then collapse the formula down to its static value
# via Copy + PasteSpecial(values), which keeps the cell's string type
# without touching its number format / style.
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row=213; A="07.05.2024"; B="07.05.2024"; C=34789; CWrap=$true;  D=14500; F=78204; I=8000;  J=498; K=124; L=4950; M="https://web.archive.org/web/20240507213334/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=214; A="08.05.2024"; B="07.05.2024"; C=34789; CWrap=$true;  D=14500; F=78204; I=8000;  J=498; K=124; L=4950; M="https://web.archive.org/web/20240508193156/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker/" },
    @{ Row=215; A="09.05.2024"; B="09.05.2024"; C=34904; CWrap=$false; D=14500; F=78514; I=8000;  J=498; K=124; L=4950; M="https://web.archive.org/web/20240509235238/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=216; A="10.05.2024"; B="10.05.2024"; C=34904; CWrap=$false; D=14500; F=78514; I=8000;  J=498; K=124; L=4950; M="https://web.archive.org/web/20240510225527/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=217; A="11.05.2024"; B="10.05.2024"; C=34904; CWrap=$false; D=14500; F=78514; I=8000;  J=498; K=124; L=4950; M="https://web.archive.org/web/20240511211611/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" },
    @{ Row=218; A="12.05.2024"; B="12.05.2024"; C=35034; CWrap=$true;  D=14500; F=78755; I=10000; J=498; K=124; L=4950; M="https://web.archive.org/web/20240512182431/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # --- Column A (tracker_date, text) ---
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Formula = "=""" + $r.A + """"
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    # --- Column B (report_date, text) ---
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Formula = "=""" + $r.B + """"
    $cellB.Copy()
    $cellB.PasteSpecial(-4163)

    # --- Column C (numeric, occasionally wrapped) ---
    $cellC = $ws.Cells.Item($row, 3)
    $cellC.Value = $r.C
    if ($r.CWrap) {
        $cellC.WrapText = $true
    }

    # --- Column D ---
    $ws.Cells.Item($row, 4).Value = $r.D

    # --- Column F (column E intentionally left blank) ---
    $ws.Cells.Item($row, 6).Value = $r.F

    # --- Columns I, J, K, L (G, H intentionally left blank) ---
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L

    # --- Column M (source url, text) ---
    $cellM = $ws.Cells.Item($row, 13)
    $cellM.Formula = "=""" + $r.M + """"
    $cellM.Copy()
    $cellM.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------------
# View-state bookkeeping, mirroring what a person would see after scrolling
# down to review/enter the newly appended rows: split the panes just above
# the new data and leave the active selection on the last populated cell.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.SplitRow = 208
$ws.Range("M218").Select() | Out-Null
